$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 0. Grab a "donor" FormattedText for a lone "." character that already has
#    the exact run properties we need later (rFonts cs=Arial, sz=16, szCs=16,
#    no eastAsia, no rsid) so we can stamp that formatting onto a freshly
#    inserted run without disturbing anything else. Do this up-front, before
#    any text in the document shifts around.
# ---------------------------------------------------------------------------
$donorOld = "esented due to assay-based allele bias). "
$donorRng = $d.Content
$donorRng.Find.Execute($donorOld, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$donorText = $donorRng.Text
$donorPeriodIdx = $donorText.IndexOf(".")
$donorPeriodStart = $donorRng.Start + $donorPeriodIdx
$donorChar = $d.Range($donorPeriodStart, $donorPeriodStart + 1)
$periodFormattedText = $donorChar.FormattedText

# ---------------------------------------------------------------------------
# 1. "Germline variant analysis of HAVCR2 ..." - two runs holding adjacent
#    text are coalesced into one (no visible text change).
# ---------------------------------------------------------------------------
$t1 = "Germline variant analysis of HAVCR2 exon 2 including Tyr82 and Ile97 hotspot variant loci."
$d.Content.Find.Execute($t1, $false, $false, $false, $false, $false, $true, 1, $false, $t1, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "Variants are analysed using PathOS software ..." - two runs coalesced.
# ---------------------------------------------------------------------------
$t2 = "Variants are analysed using PathOS software (Peter Mac) and described according to HGVS nomenclature version 19.01 (http://varnomen.hgvs.org/) with minor differences in accordance with Peter MacCallum Cancer Centre Molecular Pathology departmental policy. "
$d.Content.Find.Execute($t2, $false, $false, $false, $false, $false, $true, 1, $false, $t2, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. "(Richards et al. 2015, PMID: 25741868) ..." - two runs coalesced.
# ---------------------------------------------------------------------------
$t3 = "(Richards et al. 2015, PMID: 25741868) with class 3 (uncertain significance), class 4 (likely pathogenic) and class 5 (pathogenic) variants reported only."
$d.Content.Find.Execute($t3, $false, $false, $false, $false, $false, $true, 1, $false, $t3, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. "The detection limit of this assay ..." - many runs coalesced into one,
#    and the VAF figure is updated from 2% to 4%.
# ---------------------------------------------------------------------------
$old4 = "The detection limit of this assay for specimens sequenced to the target read depth of 250x is a variant allele frequency (VAF) of approximately 2%. This assay is primarily qualitative however, the variant read frequency (VRF) is provided to assist with variant interpretation and is assumed to approximate VAF in most instances (noting that the VAF of some insertions/deletions may be underrepresented due to assay-based allele bias). Copy number variations, loss of heterozygosity, structural rearrangements or aneuploidies are not reported. Insertions or deletions (particularly those > 25 bp in length) are not reliably detected by this assay. Genes are analysed using the reference transcripts listed below; coding exons found in alternative transcripts are not assessed by this assay. "
$new4 = "The detection limit of this assay for specimens sequenced to the target read depth of 250x is a variant allele frequency (VAF) of approximately 4%. This assay is primarily qualitative however, the variant read frequency (VRF) is provided to assist with variant interpretation and is assumed to approximate VAF in most instances (noting that the VAF of some insertions/deletions may be underrepresented due to assay-based allele bias). Copy number variations, loss of heterozygosity, structural rearrangements or aneuploidies are not reported. Insertions or deletions (particularly those > 25 bp in length) are not reliably detected by this assay. Genes are analysed using the reference transcripts listed below; coding exons found in alternative transcripts are not assessed by this assay. "
$d.Content.Find.Execute($old4, $false, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. "...variant zygosity is assumed ... for this patient" - three runs
#    coalesced into one (still no visible text change), then a trailing
#    full stop is appended in its own new run to end the sentence.
# ---------------------------------------------------------------------------
$old5 = ", variant zygosity is assumed to be either heterozygous or homozygous in the germline based on allele frequency for the purpose of clinical interpretation. Please note Peter Mac assumes sample identification, family relationships, and clinical diagnoses are as stated on the request. Our clinical recommendations may be based on evidence from third-party data sources and should be interpreted in the context of all other clinical and laboratory information for this patient"
$d.Content.Find.Execute($old5, $false, $false, $false, $false, $false, $true, 1, $false, $old5, 2) | Out-Null

$rng5 = $d.Content
$rng5.Find.Execute($old5, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$periodRng = $d.Range($rng5.End, $rng5.End)
$periodRng.InsertAfter(".")
$periodRng = $d.Range($rng5.End, $rng5.End + 1)
$periodRng.FormattedText = $periodFormattedText

# ---------------------------------------------------------------------------
# 6. "Please note variants may not be optimally detected ..." - three runs
#    coalesced into one (no visible text change).
# ---------------------------------------------------------------------------
$t6 = "Please note variants may not be optimally detected in genes with less than 100% coverage. The gene coverage above is considered acceptable given the available information about the clinical context, however please contact the laboratory for further advice should specific genes covered at less than 100% require full coverage. A list of regions with suboptimal coverage is available upon request."
$d.Content.Find.Execute($t6, $false, $false, $false, $false, $false, $true, 1, $false, $t6, 2) | Out-Null
